# The commit replaces the entire report body with a single blank
# paragraph, keeping only the section properties (page size/margins).
# The lone surviving paragraph is the pre-existing empty paragraph that
# sat between "... generated 18" and "automatically." in the original
# "Step 7" procedure text - everything else (headings, steps, the
# OUTPUT screenshot and the Result line) is removed.

$d = $word.ActiveDocument

# Find the paragraph that is already empty (its Range.Text is just the
# paragraph mark, chr(13)) - that is the one paragraph that must survive.
$emptyIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.TrimEnd([char]13).Length -eq 0) {
        $emptyIndex = $i
        break
    }
}

if ($emptyIndex -eq -1) {
    throw "Could not locate the blank paragraph to preserve"
}

# Delete everything after the blank paragraph first (so earlier offsets
# used below stay valid), then delete everything before it.
if ($emptyIndex -lt $d.Paragraphs.Count) {
    $tailStart = $d.Paragraphs.Item($emptyIndex + 1).Range.Start
    $tailEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
    $d.Range($tailStart, $tailEnd).Delete()
}

if ($emptyIndex -gt 1) {
    $headEnd = $d.Paragraphs.Item($emptyIndex - 1).Range.End
    $d.Range(0, $headEnd).Delete()
}

# Best-effort: the uploaded version also marks the built-in "Default
# Paragraph Font" character style as semi-hidden (it is no longer in use
# once the body text is gone). Harmless no-op if unsupported.
try {
    $d.Styles.Item("Default Paragraph Font").Hidden = $true
} catch {
}

Write-Output ("ParagraphsRemaining=" + $d.Paragraphs.Count)
Write-Output ("RemainingText=[" + $d.Paragraphs.Item(1).Range.Text.TrimEnd([char]13) + "]")
